# Update LR-pair TPM-derived values (Sema7a-Itga1) per new TPM recomputation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 32.21267
$ws.Range("H2").Value = 96.63801000000001
$ws.Range("I2").Value = 0.7096649552378644
$ws.Range("J2").Value = 0.7096649552378644
$ws.Range("M2").Value = 71.44418333333333
$ws.Range("N2").Value = 214.33255
$ws.Range("O2").Value = 0.6986063918429039
$ws.Range("P2").Value = 0.6986063918429037
$ws.Range("Q2").Value = 2301.407901136167
$ws.Range("R2").Value = 20712.6711102255
$ws.Range("S2").Value = 0.4957764737960804
$ws.Range("T2").Value = 0.4957764737960803
# Row 3
$ws.Range("G3").Value = 32.21267
$ws.Range("H3").Value = 96.63801000000001
$ws.Range("I3").Value = 0.7096649552378644
$ws.Range("J3").Value = 0.7096649552378644
$ws.Range("O3").Value = 0.1188372961583501
$ws.Range("P3").Value = 0.1188372961583501
$ws.Range("Q3").Value = 391.4838105145567
$ws.Range("R3").Value = 3523.35429463101
$ws.Range("S3").Value = 0.08433466445880435
$ws.Range("T3").Value = 0.08433466445880433
# Row 4
$ws.Range("G4").Value = 32.21267
$ws.Range("H4").Value = 96.63801000000001
$ws.Range("I4").Value = 0.7096649552378644
$ws.Range("J4").Value = 0.7096649552378644
$ws.Range("M4").Value = 18.571964
$ws.Range("N4").Value = 55.715892
$ws.Range("O4").Value = 0.1816032062252276
$ws.Range("P4").Value = 0.1816032062252276
$ws.Range("Q4").Value = 598.25254758388
$ws.Range("R4").Value = 5384.27292825492
$ws.Range("S4").Value = 0.1288774312168788
$ws.Range("T4").Value = 0.1288774312168788
# Row 5
$ws.Range("G5").Value = 32.21267
$ws.Range("H5").Value = 96.63801000000001
$ws.Range("I5").Value = 0.7096649552378644
$ws.Range("J5").Value = 0.7096649552378644
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.09747100000000002
$ws.Range("N5").Value = 0.292413
$ws.Range("O5").Value = 0.000953105773518577
$ws.Range("P5").Value = 0.0009531057735185768
$ws.Range("Q5").Value = 3.139801157570001
$ws.Range("R5").Value = 28.25821041813001
$ws.Range("S5").Value = 0.0006763857661010111
$ws.Range("T5").Value = 0.0006763857661010109
# Row 6
$ws.Range("I6").Value = 0.2527239295880077
$ws.Range("J6").Value = 0.2527239295880077
$ws.Range("M6").Value = 71.44418333333333
$ws.Range("N6").Value = 214.33255
$ws.Range("O6").Value = 0.6986063918429039
$ws.Range("P6").Value = 0.6986063918429037
$ws.Range("Q6").Value = 819.5710441486777
$ws.Range("R6").Value = 7376.1393973381
$ws.Range("S6").Value = 0.1765545525818382
$ws.Range("T6").Value = 0.1765545525818381
# Row 7
$ws.Range("I7").Value = 0.2527239295880077
$ws.Range("J7").Value = 0.2527239295880077
$ws.Range("O7").Value = 0.1188372961583501
$ws.Range("P7").Value = 0.1188372961583501
$ws.Range("S7").Value = 0.03003302846675208
$ws.Range("T7").Value = 0.03003302846675208
# Row 8
$ws.Range("I8").Value = 0.2527239295880077
$ws.Range("J8").Value = 0.2527239295880077
$ws.Range("M8").Value = 18.571964
$ws.Range("N8").Value = 55.715892
$ws.Range("O8").Value = 0.1816032062252276
$ws.Range("P8").Value = 0.1816032062252276
$ws.Range("Q8").Value = 213.0480497811226
$ws.Range("R8").Value = 1917.432448030104
$ws.Range("S8").Value = 0.04589547590302086
$ws.Range("T8").Value = 0.04589547590302086
# Row 9
$ws.Range("I9").Value = 0.2527239295880077
$ws.Range("J9").Value = 0.2527239295880077
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.09747100000000002
$ws.Range("N9").Value = 0.292413
$ws.Range("O9").Value = 0.000953105773518577
$ws.Range("P9").Value = 0.0009531057735185768
$ws.Range("Q9").Value = 1.118137341867333
$ws.Range("R9").Value = 10.063236076806
$ws.Range("S9").Value = 0.0002408726363966325
$ws.Range("T9").Value = 0.0002408726363966324
# Row 10
$ws.Range("G10").Value = 1.279382333333333
$ws.Range("H10").Value = 3.838147
$ws.Range("I10").Value = 0.02818558059040478
$ws.Range("J10").Value = 0.02818558059040478
$ws.Range("M10").Value = 71.44418333333333
$ws.Range("N10").Value = 214.33255
$ws.Range("O10").Value = 0.6986063918429039
$ws.Range("P10").Value = 0.6986063918429037
$ws.Range("Q10").Value = 91.40442597609443
$ws.Range("R10").Value = 822.6398337848499
$ws.Range("S10").Value = 0.01969062675826007
$ws.Range("T10").Value = 0.01969062675826007
# Row 11
$ws.Range("G11").Value = 1.279382333333333
$ws.Range("H11").Value = 3.838147
$ws.Range("I11").Value = 0.02818558059040478
$ws.Range("J11").Value = 0.02818558059040478
$ws.Range("O11").Value = 0.1188372961583501
$ws.Range("P11").Value = 0.1188372961583501
$ws.Range("Q11").Value = 15.54846186169411
$ws.Range("R11").Value = 139.936156755247
$ws.Range("S11").Value = 0.003349498188016977
$ws.Range("T11").Value = 0.003349498188016976
# Row 12
$ws.Range("G12").Value = 1.279382333333333
$ws.Range("H12").Value = 3.838147
$ws.Range("I12").Value = 0.02818558059040478
$ws.Range("J12").Value = 0.02818558059040478
$ws.Range("M12").Value = 18.571964
$ws.Range("N12").Value = 55.715892
$ws.Range("O12").Value = 0.1816032062252276
$ws.Range("P12").Value = 0.1816032062252276
$ws.Range("Q12").Value = 23.76064263690266
$ws.Range("R12").Value = 213.845783732124
$ws.Range("S12").Value = 0.005118591804537052
$ws.Range("T12").Value = 0.005118591804537051
# Row 13
$ws.Range("G13").Value = 1.279382333333333
$ws.Range("H13").Value = 3.838147
$ws.Range("I13").Value = 0.02818558059040478
$ws.Range("J13").Value = 0.02818558059040478
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.09747100000000002
$ws.Range("N13").Value = 0.292413
$ws.Range("O13").Value = 0.000953105773518577
$ws.Range("P13").Value = 0.0009531057735185768
$ws.Range("Q13").Value = 0.1247026754123333
$ws.Range("R13").Value = 1.122324078711
$ws.Range("S13").Value = 0.00002686383959068794
$ws.Range("T13").Value = 0.00002686383959068793
# Row 14
$ws.Range("G14").Value = 0.4278379999999999
$ws.Range("H14").Value = 1.283514
$ws.Range("I14").Value = 0.009425534583723031
$ws.Range("J14").Value = 0.009425534583723031
$ws.Range("M14").Value = 71.44418333333333
$ws.Range("N14").Value = 214.33255
$ws.Range("O14").Value = 0.6986063918429039
$ws.Range("P14").Value = 0.6986063918429037
$ws.Range("Q14").Value = 30.56653650896666
$ws.Range("R14").Value = 275.0988285806999
$ws.Range("S14").Value = 0.006584738706725254
$ws.Range("T14").Value = 0.006584738706725253
# Row 15
$ws.Range("G15").Value = 0.4278379999999999
$ws.Range("H15").Value = 1.283514
$ws.Range("I15").Value = 0.009425534583723031
$ws.Range("J15").Value = 0.009425534583723031
$ws.Range("O15").Value = 0.1188372961583501
$ws.Range("P15").Value = 0.1188372961583501
$ws.Range("Q15").Value = 5.199558140412666
$ws.Range("R15").Value = 46.79602326371399
$ws.Range("S15").Value = 0.001120105044776665
$ws.Range("T15").Value = 0.001120105044776665
# Row 16
$ws.Range("G16").Value = 0.4278379999999999
$ws.Range("H16").Value = 1.283514
$ws.Range("I16").Value = 0.009425534583723031
$ws.Range("J16").Value = 0.009425534583723031
$ws.Range("M16").Value = 18.571964
$ws.Range("N16").Value = 55.715892
$ws.Range("O16").Value = 0.1816032062252276
$ws.Range("P16").Value = 0.1816032062252276
$ws.Range("Q16").Value = 7.945791933831998
$ws.Range("R16").Value = 71.51212740448798
$ws.Range("S16").Value = 0.001711707300790868
$ws.Range("T16").Value = 0.001711707300790868
# Row 17
$ws.Range("G17").Value = 0.4278379999999999
$ws.Range("H17").Value = 1.283514
$ws.Range("I17").Value = 0.009425534583723031
$ws.Range("J17").Value = 0.009425534583723031
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.09747100000000002
$ws.Range("N17").Value = 0.292413
$ws.Range("O17").Value = 0.000953105773518577
$ws.Range("P17").Value = 0.0009531057735185768
$ws.Range("Q17").Value = 0.041701797698
$ws.Range("R17").Value = 0.375316179282
$ws.Range("S17").Value = 0.000008983531430245438
$ws.Range("T17").Value = 0.000008983531430245437
